$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that should advance from 46060 to 46061
# for every data row (rows 2 through 375).
$range = $ws.Range("C2:C375")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
